$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray G19 cell / row entirely
$ws.Range("G19").EntireRow.Delete() | Out-Null

# Add new label cell A8 with the value that already exists as a shared string ("alberto")
$ws.Range("A8").Value = "alberto"

# Update the selection to match the new state
$ws.Range("A9:A16").Select() | Out-Null
$ws.Range("A16").Activate() | Out-Null
